# Auto-generated Excel COM-interop script
# Applies the crypto price/volume update described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.180.47'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.785.16'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''226.20'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D8").Value = '''32.12'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '''0.292'
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").Value = '''0.0688'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '2.042.39'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''10.96'
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.754.18'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '34.176.33'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '''67.68'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").Value = '''245.86'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").Value = '0.0₃0793'
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").Value = '''10.97'
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = '''4.13'
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = '''2.05'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '''161.64'
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").Value = '''0.115'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").Value = '''0.0519'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '''3.75'
$ws.Range("E32").Value = '  +2.37%  '
$ws.Range("D33").Value = '''3.74'
$ws.Range("E33").Value = '  +3.58%  '
$ws.Range("D34").Value = '''1.79'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = '1.446.09'
$ws.Range("E35").Value = '  +3.01%  '
$ws.Range("D36").Value = '''2.60'
$ws.Range("E36").Value = '  +10.91%  '
$ws.Range("D37").Value = '''0.655'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = '''83.00'
$ws.Range("E40").Value = '  +3.66%  '
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("D43").Value = '''0.916'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = '''13.77'
$ws.Range("E44").Value = '  +3.73%  '
$ws.Range("D45").Value = '''0.0514'
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = '''1.09'
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("D48").Value = '1.942.26'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '''104.86'
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0131'
$ws.Range("E50").Value = '  -5.61%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  +0.19%  '
